# "fixing right side container"
# Insert a new column F ("Dependency_Type_Descrip" / "Business Description..")
# in front of the existing Dependency_Name / Dependency_Descrip columns, which
# shift right to G and H respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; existing F (Dependency_Name) and G
# (Dependency_Descrip) shift to G and H.
$ws.Columns("F").Insert()

# Header for the new column.
$ws.Range("F1").Value = "Dependency_Type_Descrip"

# Populate the new column for the first block of data rows (2-36) with the
# "Business Description.." placeholder text; rows 37-74 are intentionally
# left blank.
$ws.Range("F2:F36").Value = "Business Description.."

# Match the column's width to the rest of the table (as close as the host
# allows) and restore the expected selection.
$ws.Columns("F").ColumnWidth = 23.83
$null = $ws.Range("F2:F36").Select()
